# Update the "Förändrad" (C) column date values from 45178 (2023-09-09)
# to 45179 (2023-09-10) for rows 2 through 11, matching the automatic
# data-refresh commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 11; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45178) {
        $cell.Value2 = 45179
    }
}
